$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data rows to append (update through 09/09/2021 included):
# date serial, B (nuovi pos.), C (somma mobile 7gg), D (somma mobile 7gg per 100mila abitanti)
$data = @(
    @(44441, 0, 0, 0),
    @(44442, 0, 0, 0),
    @(44443, 1, 1, 26.76659528907923),
    @(44444, 0, 1, 26.76659528907923),
    @(44445, 0, 1, 26.76659528907923),
    @(44446, 0, 1, 26.76659528907923),
    @(44447, 0, 1, 26.76659528907923),
    @(44448, 0, 1, 26.76659528907923)
)

$lastRow = 366
$startRow = $lastRow + 1

for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $startRow + $i
    $vals = $data[$i]

    # Copy the formatting from the last existing row so the new row matches
    # the look (date style, borders, etc.) of the preceding data.
    $ws.Range("A$lastRow`:D$lastRow").Copy()
    $ws.Range("A$row`:D$row").PasteSpecial(-4122)  # xlPasteFormats

    $ws.Cells.Item($row, 1).Value = $vals[0]
    $ws.Cells.Item($row, 2).Value = $vals[1]
    $ws.Cells.Item($row, 3).Value = $vals[2]
    $ws.Cells.Item($row, 4).Value = $vals[3]
}

$excel.CutCopyMode = 0
